$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "500-23-12"
$ws.Range("C4").Value = "142-96-12"
$ws.Range("D4").Value = "915-87-44"
$ws.Range("E4").Value = "967-72-31"
$ws.Range("F4").Value = "944-88-25"
$ws.Range("G4").Value = "766-43-43"
